{"js": "// Replace the 25 \"a\u00f7b=\" division prompts in the worksheet table with their\n// updated values. Several source strings are not unique across the table\n// (e.g. \"79\u00f76=\" appears twice, each needing a different replacement), so\n// matches must be consumed in document order rather than via a single\n// global find/replace.\nconst replacements = [\n  [\"76\u00f75=\", \"25\u00f79=\"],\n  [\"37\u00f78=\", \"98\u00f79=\"],\n  [\"59\u00f72=\", \"56\u00f75=\"],\n  [\"79\u00f76=\", \"54\u00f78=\"],\n  [\"41\u00f77=\", \"57\u00f74=\"],\n  [\"56\u00f76=\", \"77\u00f74=\"],\n  [\"16\u00f73=\", \"69\u00f74=\"],\n  [\"97\u00f72=\", \"10\u00f75=\"],\n  [\"95\u00f74=\", \"69\u00f74=\"],\n  [\"34\u00f72=\", \"64\u00f77=\"],\n  [\"22\u00f79=\", \"27\u00f72=\"],\n  [\"66\u00f79=\", \"58\u00f74=\"],\n  [\"76\u00f73=\", \"44\u00f74=\"],\n  [\"31\u00f74=\", \"23\u00f75=\"],\n  [\"53\u00f74=\", \"52\u00f79=\"],\n  [\"46\u00f78=\", \"27\u00f77=\"],\n  [\"35\u00f77=\", \"14\u00f74=\"],\n  [\"13\u00f79=\", \"83\u00f72=\"],\n  [\"42\u00f75=\", \"46\u00f74=\"],\n  [\"40\u00f79=\", \"89\u00f78=\"],\n  [\"79\u00f76=\", \"67\u00f74=\"],\n  [\"84\u00f74=\", \"83\u00f79=\"],\n  [\"66\u00f77=\", \"69\u00f72=\"],\n  [\"71\u00f74=\", \"32\u00f75=\"],\n  [\"28\u00f79=\", \"10\u00f73=\"],\n];\n\nconst body = context.document.body;\n\n// Group the ordered replacement list by source text, keeping the order in\n// which each target value must be applied to successive matches of that\n// source text.\nconst groups = new Map();\nfor (const [src, dst] of replacements) {\n  if (!groups.has(src)) groups.set(src, []);\n  groups.get(src).push(dst);\n}\n\n// Search once per distinct source text; Word returns the hits in document\n// order, so the i-th hit gets the i-th queued replacement for that text.\nconst searchResults = new Map();\nfor (const src of groups.keys()) {\n  const results = body.search(src, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  searchResults.set(src, results);\n}\n\nawait context.sync();\n\nfor (const [src, dsts] of groups.entries()) {\n  const items = searchResults.get(src).items;\n  if (items.length !== dsts.length) {\n    throw new Error(\n      `Expected ${dsts.length} occurrence(s) of \"${src}\" but found ${items.length}.`\n    );\n  }\n  for (let i = 0; i < items.length; i++) {\n    items[i].insertText(dsts[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"a\u00f7b=\" division prompts in the worksheet table with their\n# updated values, in document order. Several source strings repeat with\n# different target values (e.g. \"79\u00f76=\" appears twice), so replacements are\n# driven off one shared Range whose Find position advances after each\n# wdReplaceOne hit -- this walks the document from top to bottom and applies\n# the queued replacement list strictly in the order the matches occur.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Start = 0\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$pairs = @(\n    ,@(\"76\u00f75=\", \"25\u00f79=\")\n    ,@(\"37\u00f78=\", \"98\u00f79=\")\n    ,@(\"59\u00f72=\", \"56\u00f75=\")\n    ,@(\"79\u00f76=\", \"54\u00f78=\")\n    ,@(\"41\u00f77=\", \"57\u00f74=\")\n    ,@(\"56\u00f76=\", \"77\u00f74=\")\n    ,@(\"16\u00f73=\", \"69\u00f74=\")\n    ,@(\"97\u00f72=\", \"10\u00f75=\")\n    ,@(\"95\u00f74=\", \"69\u00f74=\")\n    ,@(\"34\u00f72=\", \"64\u00f77=\")\n    ,@(\"22\u00f79=\", \"27\u00f72=\")\n    ,@(\"66\u00f79=\", \"58\u00f74=\")\n    ,@(\"76\u00f73=\", \"44\u00f74=\")\n    ,@(\"31\u00f74=\", \"23\u00f75=\")\n    ,@(\"53\u00f74=\", \"52\u00f79=\")\n    ,@(\"46\u00f78=\", \"27\u00f77=\")\n    ,@(\"35\u00f77=\", \"14\u00f74=\")\n    ,@(\"13\u00f79=\", \"83\u00f72=\")\n    ,@(\"42\u00f75=\", \"46\u00f74=\")\n    ,@(\"40\u00f79=\", \"89\u00f78=\")\n    ,@(\"79\u00f76=\", \"67\u00f74=\")\n    ,@(\"84\u00f74=\", \"83\u00f79=\")\n    ,@(\"66\u00f77=\", \"69\u00f72=\")\n    ,@(\"71\u00f74=\", \"32\u00f75=\")\n    ,@(\"28\u00f79=\", \"10\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $source = $pair[0]\n    $target = $pair[1]\n    $found = $rng.Find.Execute($source, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $target, $wdReplaceOne)\n    if (-not $found) {\n        throw \"Could not find occurrence of `\"$source`\" to replace with `\"$target`\".\"\n    }\n}\n\n"}
